# Localize handlebars: expose the data model (use {{data.name}} instead of {{name}})
# in the five survey-sheet prompts, then leave the workbook focused on the
# "survey" sheet with G16 selected (matching where the last edit was made).

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

$survey.Range("G2").Value = "Can {{data.name}} read and write in any language?"
$survey.Range("G3").Value = "Has {{data.name}} ever attended school?"
$survey.Range("G5").Value = "What kind(s) of school did {{data.name}} primarily attend? "
$survey.Range("G7").Value = "What is the highest grade {{data.name}} completed?"
$survey.Range("G16").Value = "Is  {{data.name}} currently in attending school?"

$survey.Activate()
$survey.Range("G16").Select()
